{"js": "// Rewrites the two right-aligned header paragraphs of the Product Backlog\n// template: paragraph 1 becomes an \"ID:\" / \"{id}\" merge field, and\n// paragraph 2 gains the \"dailyWorkTime\" / \"endDate\" merge fields plus a\n// whole new STORY card (ID / CREATION DATE / DESCRIPTION merge fields).\n\nconst OOXML_NS =\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">';\n\nfunction wrapPackage(bodyXml) {\n  return (\n    '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    OOXML_NS +\n    '<w:body>' +\n    bodyXml +\n    '</w:body></w:document>' +\n    '</pkg:xmlData></pkg:part></pkg:package>'\n  );\n}\n\n// Shared run-properties blob used for every templated merge-field value\n// (id / dailyWorkTime / endDate) in this template.\nconst FIELD_RPR =\n  '<w:rPr>' +\n  '<w:rFonts w:ascii=\"Menlo Regular\" w:hAnsi=\"Menlo Regular\" w:cs=\"Menlo Regular\"/>' +\n  '<w:color w:val=\"000000\"/>' +\n  '<w:sz w:val=\"22\"/>' +\n  '<w:szCs w:val=\"22\"/>' +\n  '</w:rPr>';\n\nconst PARA1_OOXML = wrapPackage(\n  '<w:p>' +\n    '<w:pPr><w:jc w:val=\"right\"/></w:pPr>' +\n    '<w:r><w:t>ID</w:t></w:r>' +\n    '<w:r><w:t>:</w:t></w:r>' +\n    '<w:r><w:br/><w:t>{</w:t></w:r>' +\n    '<w:r>' + FIELD_RPR + '<w:t>id</w:t></w:r>' +\n    '<w:r><w:t>}</w:t></w:r>' +\n  '</w:p>'\n);\n\nconst PARA2_OOXML = wrapPackage(\n  '<w:p>' +\n    '<w:pPr><w:jc w:val=\"right\"/></w:pPr>' +\n    '<w:r><w:rPr><w:b/></w:rPr><w:t>EQUIPE</w:t></w:r>' +\n    '<w:r><w:br/><w:t>{</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r>' + FIELD_RPR + '<w:t>dailyWorkTime</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t>}</w:t></w:r>' +\n    '<w:r><w:br/></w:r>' +\n    '<w:r><w:br/></w:r>' +\n    '<w:r><w:t xml:space=\"preserve\">Nome do </w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:t>projeto</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:br/></w:r>' +\n    '<w:r><w:t>{</w:t></w:r>' +\n    '<w:r>' + FIELD_RPR + '<w:t>endDate</w:t></w:r>' +\n    '<w:r><w:t>}</w:t></w:r>' +\n    '<w:r><w:br/></w:r>' +\n    '<w:r><w:br/></w:r>' +\n    '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>' +\n    '<w:r><w:rPr><w:b/><w:sz w:val=\"32\"/><w:szCs w:val=\"32\"/></w:rPr><w:t>STORY</w:t></w:r>' +\n    '<w:r><w:br/></w:r>' +\n    '<w:r><w:t>ID</w:t></w:r>' +\n    '<w:r><w:br/><w:t>{</w:t></w:r>' +\n    '<w:r><w:t>id</w:t></w:r>' +\n    '<w:r><w:t>}</w:t></w:r>' +\n    '<w:r><w:br/></w:r>' +\n    '<w:r><w:br/><w:t>CREATION DATE</w:t></w:r>' +\n    '<w:r><w:br/><w:t>{</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellStart\"/>' +\n    '<w:r><w:t>creationDate</w:t></w:r>' +\n    '<w:proofErr w:type=\"spellEnd\"/>' +\n    '<w:r><w:t>}</w:t></w:r>' +\n    '<w:r><w:br/></w:r>' +\n    '<w:r><w:br/><w:t>DESCRIPTION</w:t></w:r>' +\n    '<w:r><w:br/><w:t>{</w:t></w:r>' +\n    '<w:r><w:t>description</w:t></w:r>' +\n    '<w:r><w:t>}</w:t></w:r>' +\n    '<w:bookmarkEnd w:id=\"0\"/>' +\n  '</w:p>'\n);\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst para1 = paragraphs.items[0]; // \"Author:\" / \"{Autor}\"\nconst para2 = paragraphs.items[1]; // \"EQUIPE\" / \"{Equipe}\" / \"Nome do projeto\" / \"{NomeProjeto}\"\n\npara1.insertOoxml(PARA1_OOXML, Word.InsertLocation.replace);\npara2.insertOoxml(PARA2_OOXML, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Rewrites the two right-aligned header paragraphs of the Product Backlog\n# template: paragraph 1 becomes an \"ID:\" / \"{id}\" merge field, and\n# paragraph 2 gains the \"dailyWorkTime\" / \"endDate\" merge fields plus a\n# whole new STORY card (ID / CREATION DATE / DESCRIPTION merge fields).\n\n$d = $word.ActiveDocument\n\n# Shared run-properties blob used for every templated merge-field value\n# (id / dailyWorkTime / endDate) in this template.\n$fieldRPr = '<w:rPr><w:rFonts w:ascii=\"Menlo Regular\" w:hAnsi=\"Menlo Regular\" w:cs=\"Menlo Regular\"/><w:color w:val=\"000000\"/><w:sz w:val=\"22\"/><w:szCs w:val=\"22\"/></w:rPr>'\n\n$para1Body = '<w:p><w:pPr><w:jc w:val=\"right\"/></w:pPr>' +\n  '<w:r><w:t>ID</w:t></w:r>' +\n  '<w:r><w:t>:</w:t></w:r>' +\n  '<w:r><w:br/><w:t>{</w:t></w:r>' +\n  '<w:r>' + $fieldRPr + '<w:t>id</w:t></w:r>' +\n  '<w:r><w:t>}</w:t></w:r>' +\n  '</w:p>'\n\n$para2Body = '<w:p><w:pPr><w:jc w:val=\"right\"/></w:pPr>' +\n  '<w:r><w:rPr><w:b/></w:rPr><w:t>EQUIPE</w:t></w:r>' +\n  '<w:r><w:br/><w:t>{</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  '<w:r>' + $fieldRPr + '<w:t>dailyWorkTime</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:t>}</w:t></w:r>' +\n  '<w:r><w:br/></w:r>' +\n  '<w:r><w:br/></w:r>' +\n  '<w:r><w:t xml:space=\"preserve\">Nome do </w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  '<w:r><w:t>projeto</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:br/></w:r>' +\n  '<w:r><w:t>{</w:t></w:r>' +\n  '<w:r>' + $fieldRPr + '<w:t>endDate</w:t></w:r>' +\n  '<w:r><w:t>}</w:t></w:r>' +\n  '<w:r><w:br/></w:r>' +\n  '<w:r><w:br/></w:r>' +\n  '<w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>' +\n  '<w:r><w:rPr><w:b/><w:sz w:val=\"32\"/><w:szCs w:val=\"32\"/></w:rPr><w:t>STORY</w:t></w:r>' +\n  '<w:r><w:br/></w:r>' +\n  '<w:r><w:t>ID</w:t></w:r>' +\n  '<w:r><w:br/><w:t>{</w:t></w:r>' +\n  '<w:r><w:t>id</w:t></w:r>' +\n  '<w:r><w:t>}</w:t></w:r>' +\n  '<w:r><w:br/></w:r>' +\n  '<w:r><w:br/><w:t>CREATION DATE</w:t></w:r>' +\n  '<w:r><w:br/><w:t>{</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellStart\"/>' +\n  '<w:r><w:t>creationDate</w:t></w:r>' +\n  '<w:proofErr w:type=\"spellEnd\"/>' +\n  '<w:r><w:t>}</w:t></w:r>' +\n  '<w:r><w:br/></w:r>' +\n  '<w:r><w:br/><w:t>DESCRIPTION</w:t></w:r>' +\n  '<w:r><w:br/><w:t>{</w:t></w:r>' +\n  '<w:r><w:t>description</w:t></w:r>' +\n  '<w:r><w:t>}</w:t></w:r>' +\n  '<w:bookmarkEnd w:id=\"0\"/>' +\n  '</w:p>'\n\nfunction Wrap-Ooxml($bodyXml) {\n  return '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData>' +\n    '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n    '<w:body>' + $bodyXml + '</w:body></w:document>' +\n    '</pkg:xmlData></pkg:part></pkg:package>'\n}\n\n$para1Range = $d.Paragraphs(1).Range\n$para1Range.InsertXML((Wrap-Ooxml $para1Body))\n\n$para2Range = $d.Paragraphs(2).Range\n$para2Range.InsertXML((Wrap-Ooxml $para2Body))\n"}
